# Applies the cryptos-list refresh described in the commit diff.
# Column layout: A=index(unchanged) B=Coin C=Link D=Price E=Volume(1h)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Cells.Item(2, 4).Value = "25.827.39"
$ws.Cells.Item(2, 5).Value = "  -0.22%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.640.59"
$ws.Cells.Item(3, 5).Value = "  +0.55%  "

# Row 4
$ws.Cells.Item(4, 4).Formula = "'1.003"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  +0.18%  "

# Row 5
$ws.Cells.Item(5, 5).Value = "  -0.03%  "

# Row 6
$ws.Cells.Item(6, 4).Formula = "'0.5071"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -0.59%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  +0.10%  "

# Row 8
$ws.Cells.Item(8, 4).Formula = "'0.2585"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +0.72%  "

# Row 9
$ws.Cells.Item(9, 4).Formula = "'0.06434"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +1.57%  "

# Row 10
$ws.Cells.Item(10, 5).Value = "  +5.39%  "

# Row 11
$ws.Cells.Item(11, 4).Formula = "'0.07796"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +0.26%  "

# Row 12
$ws.Cells.Item(12, 4).Formula = "'4.268"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +0.10%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "1.645.63"
$ws.Cells.Item(13, 5).Value = "  +0.61%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "1.866.86"
$ws.Cells.Item(14, 5).Value = "  +0.52%  "

# Row 15
$ws.Cells.Item(15, 4).Formula = "'0.5626"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +2.37%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "0.0₅7663"
$ws.Cells.Item(16, 5).Value = "  +0.40%  "

# Row 17
$ws.Cells.Item(17, 5).Value = "  -0.52%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "25.835.37"
$ws.Cells.Item(18, 5).Value = "  -0.33%  "

# Row 19
$ws.Cells.Item(19, 5).Value = "  +0.17%  "

# Row 20
$ws.Cells.Item(20, 2).Value = "Uniswap"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Cells.Item(20, 4).Formula = "'4.388"
$ws.Cells.Item(20, 4).Style = "Normal"

# Row 21
$ws.Cells.Item(21, 2).Value = "BitcoinCash"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Cells.Item(21, 4).Formula = "'193.07"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -0.57%  "

# Row 22
$ws.Cells.Item(22, 4).Formula = "'9.943"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +0.96%  "

# Row 23
$ws.Cells.Item(23, 4).Formula = "'6.149"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +2.05%  "

# Row 24
$ws.Cells.Item(24, 4).Formula = "'1.003"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +0.03%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  -4.77%  "

# Row 26
$ws.Cells.Item(26, 4).Formula = "'139.40"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -1.94%  "

# Row 27
$ws.Cells.Item(27, 4).Formula = "'0.1237"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -1.00%  "

# Row 28
$ws.Cells.Item(28, 4).Formula = "'6.836"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +1.34%  "

# Row 29
$ws.Cells.Item(29, 4).Formula = "'15.59"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +0.48%  "

# Row 30
$ws.Cells.Item(30, 5).Value = "  +0.28%  "

# Row 31
$ws.Cells.Item(31, 4).Formula = "'0.04965"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +1.69%  "

# Row 32
$ws.Cells.Item(32, 4).Formula = "'3.295"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +1.99%  "

# Row 33
$ws.Cells.Item(33, 5).Value = "  +2.19%  "

# Row 34
$ws.Cells.Item(34, 5).Value = "  +2.03%  "

# Row 35
$ws.Cells.Item(35, 4).Formula = "'2.384"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +0.38%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  +0.77%  "

# Row 37
$ws.Cells.Item(37, 4).Formula = "'2.572"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +1.31%  "

# Row 38
$ws.Cells.Item(38, 5).Value = "  +1.37%  "

# Row 39
$ws.Cells.Item(39, 4).Value = "1.134.48"
$ws.Cells.Item(39, 5).Value = "  +1.73%  "

# Row 40
$ws.Cells.Item(40, 4).Formula = "'0.01574"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +1.27%  "

# Row 41
$ws.Cells.Item(41, 4).Formula = "'0.9963"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -0.43%  "

# Row 42
$ws.Cells.Item(42, 4).Formula = "'5.480"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -1.70%  "

# Row 43
$ws.Cells.Item(43, 4).Formula = "'0.8027"
$ws.Cells.Item(43, 4).Style = "Normal"

# Row 44
$ws.Cells.Item(44, 4).Formula = "'99.03"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +1.56%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "1.777.23"
$ws.Cells.Item(45, 5).Value = "  +0.62%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "0.0₈111"
$ws.Cells.Item(46, 5).Value = "  -4.97%  "

# Row 47
$ws.Cells.Item(47, 4).Formula = "'55.65"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +1.82%  "

# Row 48
$ws.Cells.Item(48, 4).Formula = "'0.4268"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -3.83%  "

# Row 49
$ws.Cells.Item(49, 4).Formula = "'7.807"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +3.73%  "

# Row 50
$ws.Cells.Item(50, 4).Formula = "'0.05037"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -1.82%  "

# Row 51
$ws.Cells.Item(51, 4).Formula = "'0.9956"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -0.60%  "

